# The presentation's Design/theme (ppt/theme/theme1.xml, used by the one
# slide master - "Integral") is swapped for the stock "Office Theme"
# palette (previously sitting unused in ppt/theme/theme2.xml, only wired
# to the notes master). The font scheme and format scheme of the two
# themes are already identical, so only the 12-slot theme colour scheme
# needs to change.
#
# Helper: turn an "RRGGBB" hex string into the integer VBA's ColorFormat.RGB
# setter expects (R + G*256 + B*65536).
function HexToVbRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

# Office Theme colour scheme, in clrScheme slot order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeTheme = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToVbRgb($officeTheme[$i - 1])
}
